$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log row appended by the 2025-11-02 03:58:56 UTC update.
$newRow = 7

$ws.Cells.Item($newRow, 1).Value = "2025-11-02 03:58:56"

# B7 ("2025-10-30") would be auto-parsed as a date serial by the Value
# setter, so force text interpretation via NumberFormat, then restore the
# default "Normal" style so the cell ends up unstyled (same as the rest
# of the data rows) while keeping the literal string value.
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).Value = "2025-10-30"
$ws.Cells.Item($newRow, 2).Style = "Normal"

$ws.Cells.Item($newRow, 3).Value = "https://rashtriyametal.com/wp-content/uploads/2025/11/ListPrice30102025.pdf"
$ws.Cells.Item($newRow, 4).Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/ListPrice30102025.pdf"
